# Sports Info Solutions Analytics Challenge.pptx -- "Add files via upload"
#
# 1) The deck's dynamic "datetimeFigureOut" date field (Insert > Header &
#    Footer > Date) was re-cached from 7/20/2021 to 7/21/2021. That field
#    lives on the slide master and on every slide layout (12 copies total).
# 2) On slide 16 ("Our Best Route Combos"), four shapes were repositioned /
#    resized to make room (title box, the two route-combo chart pictures,
#    and the content placeholder with the write-up).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update every "Date Placeholder" (ppPlaceholderDate = 16) field's
#    cached text, on the slide master and on each of its custom layouts.
# ---------------------------------------------------------------------
function Update-DatePlaceholders($container, [string]$newText) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $shp.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master "7/21/2021"

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout "7/21/2021"
}

# ---------------------------------------------------------------------
# 2) Slide 16: reflow the title, the two pictures, and the content
#    placeholder to the new layout.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(16)

# "Title 1" -- shrink / slide right to make room for the enlarged photos.
$title = $slide.Shapes.Item(2)
$title.Left = 655.9266141732284
$title.Top = 49.99574803149606
$title.Width = 253.50188976377953
$title.Height = 114.23283464566929

# "Picture 13" -- the left route-combo chart, enlarged to fill the slide.
$pic13 = $slide.Shapes.Item(3)
$pic13.Left = 15.728425196850393
$pic13.Top = 11.660551181102361
$pic13.Width = 334.2674803149606
$pic13.Height = 484.44582677165357

# "Picture 9" -- the middle route-combo chart, enlarged to match.
$pic9 = $slide.Shapes.Item(5)
$pic9.Left = 359.13937007874017
$pic9.Top = 7.509055118110236
$pic9.Width = 278.5003937007874
$pic9.Height = 488.5973228346457

# "Content Placeholder 2" -- the write-up text box, narrowed / shifted right.
$content = $slide.Shapes.Item(6)
$content.Left = 663.8532283464567
$content.Top = 173.1428346456693
$content.Width = 245.57527559055117
$content.Height = 288.9905511811024
